# Regenerate the handback-status report: drop the completed
# "6100965e-3277-4e74-8ceb-b89abe4613f0" row from every sheet and refresh
# the handback timestamps for the file that is still pending.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A3").Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

# --- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E2").Value = "2016-03-19 16:48:14"
$ws.Range("H2").Value = "2016-03-19 16:48:59"
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()
$ws.Range("F3").Hyperlinks.Delete()
$ws.Range("G3").Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("E2").Value = "2016-03-19 16:48:24"
$ws.Range("H2").Value = "2016-03-19 16:49:14"
$ws.Range("A3").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()
$ws.Range("F3").Hyperlinks.Delete()
$ws.Range("G3").Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()
